# Daily attendance processing - 2025-10-01 19:15:12
# Applies the day's attendance-system updates to the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Recorded By" annotations: prefix/merge in the attendance-system account ---
$ws.Range("G2").Value  = "system, backup@backdoor.com"
$ws.Range("G6").Value  = "System, dnasr281@gmail.com"
$ws.Range("G12").Value = "System, dnasr281@gmail.com"
$ws.Range("G13").Value = "System, dnasr281@gmail.com"
$ws.Range("G29").Value = "system, backup@backdoor.com"
$ws.Range("G33").Value = "System, dnasr281@gmail.com"
$ws.Range("G39").Value = "System, dnasr281@gmail.com"
$ws.Range("G40").Value = "System, dnasr281@gmail.com"
$ws.Range("G56").Value = "system, backup@backdoor.com"
$ws.Range("G60").Value = "System, dnasr281@gmail.com"
$ws.Range("G66").Value = "System, dnasr281@gmail.com"
$ws.Range("G67").Value = "System, dnasr281@gmail.com"

# Reorder multi-recorder lists so admin@admin.com sorts first
$ws.Range("G90").Value  = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G116").Value = "admin@admin.com, dnasr281@gmail.com"
$ws.Range("G142").Value = "admin@admin.com, dnasr281@gmail.com"

# --- Updated attendance counts on previously recorded sessions ---
$ws.Range("H2").Value  = "33/53"
$ws.Range("H6").Value  = "44/53"
$ws.Range("H12").Value = "31/53"
$ws.Range("H13").Value = "36/53"
$ws.Range("H92").Value  = "43/56"
$ws.Range("H118").Value = "45/55"

# --- Class Statistics block (K/L columns) ---
$ws.Range("L6").Value = 68
$ws.Range("L7").Value = 1
$ws.Range("L9").Value  = "42.8%"
$ws.Range("L10").Value = "62.2%"

# --- Per-student summary table (M:S columns) ---
$ws.Range("S15").Value = "60.1%"

$ws.Range("O18").Value = 11
$ws.Range("P18").Value = 0
$ws.Range("R18").Value = "42.3%"
$ws.Range("S18").Value = "64.4%"

$ws.Range("O19").Value = 11
$ws.Range("P19").Value = 0
$ws.Range("R19").Value = "42.3%"
$ws.Range("S19").Value = "67.6%"

$ws.Range("O20").Value = 11
$ws.Range("P20").Value = 0
$ws.Range("R20").Value = "42.3%"
$ws.Range("S20").Value = "71.6%"

# --- Newly-recorded sessions: rows 93, 119 and 145 flip from "Not Recorded" ---
# (pink highlight) to the normal "Recorded" look, and gain recorder/attendance data.

$ws.Range("A92:I92").Copy()
$ws.Range("A93:I93").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G93").Value = "dnasr281@gmail.com"
$ws.Range("H93").Value = "46/56"
$ws.Range("I93").Value = "Recorded"

$ws.Range("A118:I118").Copy()
$ws.Range("A119:I119").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G119").Value = "dnasr281@gmail.com"
$ws.Range("H119").Value = "40/55"
$ws.Range("I119").Value = "Recorded"

$ws.Range("A144:I144").Copy()
$ws.Range("A145:I145").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G145").Value = "dnasr281@gmail.com"
$ws.Range("H145").Value = "48/57"
$ws.Range("I145").Value = "Recorded"
